$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 10.27784
$ws.Range("H2").Value = 30.83352
$ws.Range("I2").Value = 0.230301226653591
$ws.Range("J2").Value = 0.230301226653591
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 8.452892666666665
$ws.Range("N2").Value = 25.358678
$ws.Range("O2").Value = 0.5664982795292011
$ws.Range("P2").Value = 0.566498279529201
$ws.Range("Q2").Value = 86.87747836517332
$ws.Range("R2").Value = 781.8973052865599
$ws.Range("S2").Value = 0.1304652486727239
$ws.Range("T2").Value = 0.1304652486727239

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 10.27784
$ws.Range("H3").Value = 30.83352
$ws.Range("I3").Value = 0.230301226653591
$ws.Range("J3").Value = 0.230301226653591
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 3.429517666666667
$ws.Range("N3").Value = 10.288553
$ws.Range("O3").Value = 0.2298403557687432
$ws.Range("P3").Value = 0.2298403557687431
$ws.Range("Q3").Value = 35.24803385517333
$ws.Range("R3").Value = 317.23230469656
$ws.Range("S3").Value = 0.05293251586803932
$ws.Range("T3").Value = 0.05293251586803931

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 10.27784
$ws.Range("H4").Value = 30.83352
$ws.Range("I4").Value = 0.230301226653591
$ws.Range("J4").Value = 0.230301226653591
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.6234873333333334
$ws.Range("N4").Value = 1.870462
$ws.Range("O4").Value = 0.04178504514015868
$ws.Range("P4").Value = 0.04178504514015867
$ws.Range("Q4").Value = 6.408103054026667
$ws.Range("R4").Value = 57.67292748624
$ws.Range("S4").Value = 0.009623147151554216
$ws.Range("T4").Value = 0.009623147151554212

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 10.27784
$ws.Range("H5").Value = 30.83352
$ws.Range("I5").Value = 0.230301226653591
$ws.Range("J5").Value = 0.230301226653591
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 2.415405666666667
$ws.Range("N5").Value = 7.246217000000001
$ws.Range("O5").Value = 0.1618763195618971
$ws.Range("P5").Value = 0.1618763195618971
$ws.Range("Q5").Value = 24.82515297709334
$ws.Range("R5").Value = 223.42637679384
$ws.Range("S5").Value = 0.03728031496127359
$ws.Range("T5").Value = 0.03728031496127358

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 14.56812733333333
$ws.Range("H6").Value = 43.704382
$ws.Range("I6").Value = 0.3264360600001921
$ws.Range("J6").Value = 0.326436060000192
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 8.452892666666665
$ws.Range("N6").Value = 25.358678
$ws.Range("O6").Value = 0.5664982795292011
$ws.Range("P6").Value = 0.566498279529201
$ws.Range("Q6").Value = 123.1428167029995
$ws.Range("R6").Value = 1108.285350326996
$ws.Range("S6").Value = 0.1849254663663999
$ws.Range("T6").Value = 0.1849254663663998

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 14.56812733333333
$ws.Range("H7").Value = 43.704382
$ws.Range("I7").Value = 0.3264360600001921
$ws.Range("J7").Value = 0.326436060000192
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 3.429517666666667
$ws.Range("N7").Value = 10.288553
$ws.Range("O7").Value = 0.2298403557687432
$ws.Range("P7").Value = 0.2298403557687431
$ws.Range("Q7").Value = 49.96165005991622
$ws.Range("R7").Value = 449.654850539246
$ws.Range("S7").Value = 0.07502818016619095
$ws.Range("T7").Value = 0.07502818016619092

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 14.56812733333333
$ws.Range("H8").Value = 43.704382
$ws.Range("I8").Value = 0.3264360600001921
$ws.Range("J8").Value = 0.326436060000192
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 0.6234873333333334
$ws.Range("N8").Value = 1.870462
$ws.Range("O8").Value = 0.04178504514015868
$ws.Range("P8").Value = 0.04178504514015867
$ws.Range("Q8").Value = 9.083042862720443
$ws.Range("R8").Value = 81.74738576448399
$ws.Range("S8").Value = 0.01364014550248357
$ws.Range("T8").Value = 0.01364014550248357

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 14.56812733333333
$ws.Range("H9").Value = 43.704382
$ws.Range("I9").Value = 0.3264360600001921
$ws.Range("J9").Value = 0.326436060000192
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 2.415405666666667
$ws.Range("N9").Value = 7.246217000000001
$ws.Range("O9").Value = 0.1618763195618971
$ws.Range("P9").Value = 0.1618763195618971
$ws.Range("Q9").Value = 35.18793731365489
$ws.Range("R9").Value = 316.691435822894
$ws.Range("S9").Value = 0.05284226796511771
$ws.Range("T9").Value = 0.05284226796511769

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 2.650137
$ws.Range("H10").Value = 7.950411
$ws.Range("I10").Value = 0.05938308067649115
$ws.Range("J10").Value = 0.05938308067649114
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 8.452892666666665
$ws.Range("N10").Value = 25.358678
$ws.Range("O10").Value = 0.5664982795292011
$ws.Range("P10").Value = 0.566498279529201
$ws.Range("Q10").Value = 22.401323612962
$ws.Range("R10").Value = 201.611912516658
$ws.Range("S10").Value = 0.03364041303637598
$ws.Range("T10").Value = 0.03364041303637597

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 2.650137
$ws.Range("H11").Value = 7.950411
$ws.Range("I11").Value = 0.05938308067649115
$ws.Range("J11").Value = 0.05938308067649114
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 3.429517666666667
$ws.Range("N11").Value = 10.288553
$ws.Range("O11").Value = 0.2298403557687432
$ws.Range("P11").Value = 0.2298403557687431
$ws.Range("Q11").Value = 9.088691660587001
$ws.Range("R11").Value = 81.798224945283
$ws.Range("S11").Value = 0.0136486283893287
$ws.Range("T11").Value = 0.0136486283893287

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 2.650137
$ws.Range("H12").Value = 7.950411
$ws.Range("I12").Value = 0.05938308067649115
$ws.Range("J12").Value = 0.05938308067649114
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 0.6234873333333334
$ws.Range("N12").Value = 1.870462
$ws.Range("O12").Value = 0.04178504514015868
$ws.Range("P12").Value = 0.04178504514015867
$ws.Range("Q12").Value = 1.652326851098
$ws.Range("R12").Value = 14.870941659882
$ws.Range("S12").Value = 0.002481324706628867
$ws.Range("T12").Value = 0.002481324706628866

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 2.650137
$ws.Range("H13").Value = 7.950411
$ws.Range("I13").Value = 0.05938308067649115
$ws.Range("J13").Value = 0.05938308067649114
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 2.415405666666667
$ws.Range("N13").Value = 7.246217000000001
$ws.Range("O13").Value = 0.1618763195618971
$ws.Range("P13").Value = 0.1618763195618971
$ws.Range("Q13").Value = 6.401155927243001
$ws.Range("R13").Value = 57.610403345187
$ws.Range("S13").Value = 0.009612714544157597
$ws.Range("T13").Value = 0.009612714544157595

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 17.13170833333333
$ws.Range("H14").Value = 51.395125
$ws.Range("I14").Value = 0.3838796326697257
$ws.Range("J14").Value = 0.3838796326697257
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 8.452892666666665
$ws.Range("N14").Value = 25.358678
$ws.Range("O14").Value = 0.5664982795292011
$ws.Range("P14").Value = 0.566498279529201
$ws.Range("Q14").Value = 144.8124917383055
$ws.Range("R14").Value = 1303.31242564475
$ws.Range("S14").Value = 0.2174671514537013
$ws.Range("T14").Value = 0.2174671514537012

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 17.13170833333333
$ws.Range("H15").Value = 51.395125
$ws.Range("I15").Value = 0.3838796326697257
$ws.Range("J15").Value = 0.3838796326697257
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 3.429517666666667
$ws.Range("N15").Value = 10.288553
$ws.Range("O15").Value = 0.2298403557687432
$ws.Range("P15").Value = 0.2298403557687431
$ws.Range("Q15").Value = 58.75349638934722
$ws.Range("R15").Value = 528.781467504125
$ws.Range("S15").Value = 0.08823103134518422
$ws.Range("T15").Value = 0.08823103134518419

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 17.13170833333333
$ws.Range("H16").Value = 51.395125
$ws.Range("I16").Value = 0.3838796326697257
$ws.Range("J16").Value = 0.3838796326697257
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 0.6234873333333334
$ws.Range("N16").Value = 1.870462
$ws.Range("O16").Value = 0.04178504514015868
$ws.Range("P16").Value = 0.04178504514015867
$ws.Range("Q16").Value = 10.68140314419444
$ws.Range("R16").Value = 96.13262829775
$ws.Range("S16").Value = 0.01604042777949202
$ws.Range("T16").Value = 0.01604042777949202

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 17.13170833333333
$ws.Range("H17").Value = 51.395125
$ws.Range("I17").Value = 0.3838796326697257
$ws.Range("J17").Value = 0.3838796326697257
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 2.415405666666667
$ws.Range("N17").Value = 7.246217000000001
$ws.Range("O17").Value = 0.1618763195618971
$ws.Range("P17").Value = 0.1618763195618971
$ws.Range("Q17").Value = 41.38002538801389
$ws.Range("R17").Value = 372.420228492125
$ws.Range("S17").Value = 0.0621410220913482
$ws.Range("T17").Value = 0.0621410220913482
